{"js": "// Update the date heading in the first paragraph of the document body.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nconst titlePara = paragraphs.items[0];\nconst titleResults = titlePara.search(\"2023-11-03 Friday\", {matchWholeWord: false});\ntitleResults.load(\"items\");\nawait context.sync();\ntitleResults.items[0].insertText(\"2023-11-04 Saturday\", Word.InsertLocation.replace);\n\n// Update each division-problem cell in the table, addressed by (row, column).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst cellEdits = [\n  { row: 0, col: 0, oldText: \"67\u00f75=13, 2\", newText: \"71\u00f74=17, 3\" },\n  { row: 0, col: 1, oldText: \"18\u00f73=6, 0\", newText: \"70\u00f78=8, 6\" },\n  { row: 0, col: 2, oldText: \"52\u00f75=10, 2\", newText: \"72\u00f74=18, 0\" },\n  { row: 0, col: 3, oldText: \"30\u00f73=10, 0\", newText: \"36\u00f79=4, 0\" },\n  { row: 0, col: 4, oldText: \"87\u00f77=12, 3\", newText: \"52\u00f77=7, 3\" },\n  { row: 4, col: 0, oldText: \"16\u00f72=8, 0\", newText: \"38\u00f79=4, 2\" },\n  { row: 4, col: 1, oldText: \"86\u00f73=28, 2\", newText: \"47\u00f78=5, 7\" },\n  { row: 4, col: 2, oldText: \"87\u00f73=29, 0\", newText: \"81\u00f79=9, 0\" },\n  { row: 4, col: 3, oldText: \"57\u00f78=7, 1\", newText: \"29\u00f79=3, 2\" },\n  { row: 4, col: 4, oldText: \"76\u00f74=19, 0\", newText: \"86\u00f76=14, 2\" },\n  { row: 8, col: 0, oldText: \"17\u00f75=3, 2\", newText: \"60\u00f79=6, 6\" },\n  { row: 8, col: 1, oldText: \"90\u00f78=11, 2\", newText: \"99\u00f73=33, 0\" },\n  { row: 8, col: 2, oldText: \"59\u00f74=14, 3\", newText: \"55\u00f77=7, 6\" },\n  { row: 8, col: 3, oldText: \"75\u00f73=25, 0\", newText: \"54\u00f74=13, 2\" },\n  { row: 8, col: 4, oldText: \"13\u00f73=4, 1\", newText: \"30\u00f72=15, 0\" },\n  { row: 12, col: 0, oldText: \"79\u00f79=8, 7\", newText: \"48\u00f77=6, 6\" },\n  { row: 12, col: 1, oldText: \"96\u00f73=32, 0\", newText: \"80\u00f72=40, 0\" },\n  { row: 12, col: 2, oldText: \"64\u00f74=16, 0\", newText: \"38\u00f77=5, 3\" },\n  { row: 12, col: 3, oldText: \"72\u00f76=12, 0\", newText: \"12\u00f73=4, 0\" },\n  { row: 12, col: 4, oldText: \"42\u00f79=4, 6\", newText: \"15\u00f78=1, 7\" },\n  { row: 16, col: 0, oldText: \"39\u00f73=13, 0\", newText: \"30\u00f77=4, 2\" },\n  { row: 16, col: 1, oldText: \"81\u00f79=9, 0\", newText: \"67\u00f77=9, 4\" },\n  { row: 16, col: 2, oldText: \"24\u00f74=6, 0\", newText: \"17\u00f78=2, 1\" },\n  { row: 16, col: 3, oldText: \"18\u00f79=2, 0\", newText: \"32\u00f72=16, 0\" },\n  { row: 16, col: 4, oldText: \"22\u00f75=4, 2\", newText: \"39\u00f76=6, 3\" },\n];\n\nfor (const edit of cellEdits) {\n  const cell = table.getCell(edit.row, edit.col);\n  const results = cell.body.search(edit.oldText, { matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(edit.newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading in the first paragraph\n$d.Paragraphs.Item(1).Range.Find.Execute(\"2023-11-03 Friday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2023-11-04 Saturday\", 2) | Out-Null\n\n# Update each division-problem cell in the table, addressed by (row, column)\n$t = $d.Tables.Item(1)\n$t.Cell(1, 1).Range.Text = \"71\u00f74=17, 3\"\n$t.Cell(1, 2).Range.Text = \"70\u00f78=8, 6\"\n$t.Cell(1, 3).Range.Text = \"72\u00f74=18, 0\"\n$t.Cell(1, 4).Range.Text = \"36\u00f79=4, 0\"\n$t.Cell(1, 5).Range.Text = \"52\u00f77=7, 3\"\n$t.Cell(5, 1).Range.Text = \"38\u00f79=4, 2\"\n$t.Cell(5, 2).Range.Text = \"47\u00f78=5, 7\"\n$t.Cell(5, 3).Range.Text = \"81\u00f79=9, 0\"\n$t.Cell(5, 4).Range.Text = \"29\u00f79=3, 2\"\n$t.Cell(5, 5).Range.Text = \"86\u00f76=14, 2\"\n$t.Cell(9, 1).Range.Text = \"60\u00f79=6, 6\"\n$t.Cell(9, 2).Range.Text = \"99\u00f73=33, 0\"\n$t.Cell(9, 3).Range.Text = \"55\u00f77=7, 6\"\n$t.Cell(9, 4).Range.Text = \"54\u00f74=13, 2\"\n$t.Cell(9, 5).Range.Text = \"30\u00f72=15, 0\"\n$t.Cell(13, 1).Range.Text = \"48\u00f77=6, 6\"\n$t.Cell(13, 2).Range.Text = \"80\u00f72=40, 0\"\n$t.Cell(13, 3).Range.Text = \"38\u00f77=5, 3\"\n$t.Cell(13, 4).Range.Text = \"12\u00f73=4, 0\"\n$t.Cell(13, 5).Range.Text = \"15\u00f78=1, 7\"\n$t.Cell(17, 1).Range.Text = \"30\u00f77=4, 2\"\n$t.Cell(17, 2).Range.Text = \"67\u00f77=9, 4\"\n$t.Cell(17, 3).Range.Text = \"17\u00f78=2, 1\"\n$t.Cell(17, 4).Range.Text = \"32\u00f72=16, 0\"\n$t.Cell(17, 5).Range.Text = \"39\u00f76=6, 3\"\n"}
